# Generate Report for Handback
# Update handback status timestamps and priority value to reflect a
# fresh report generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G)
$wsOverview.Range("G2").Value = "2016-09-03 06:18:31"
$wsOverview.Range("G5").Value = "2016-09-03 06:18:31"

# zh-cn sheet: Priority (E) ht -> mt, Handoff/Handback datetimes (H, K)
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-09-03 06:18:26"
$wsZhCn.Range("H5").Value = "2016-09-03 06:18:26"
$wsZhCn.Range("K2").Value = "2016-09-03 06:18:42"
$wsZhCn.Range("K5").Value = "2016-09-03 06:18:42"

# de-de sheet: Priority (E) ht -> mt, Handoff datetime (H), Handback datetime (K)
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-09-03 06:18:31"
$wsDeDe.Range("H5").Value = "2016-09-03 06:18:31"
$wsDeDe.Range("K2").Value = "2016-09-03 06:18:49"
$wsDeDe.Range("K5").Value = "2016-09-03 06:18:49"
